$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: duplicate row 21's current formatting (style 4/5) down onto new row 22 ---
$ws.Range("A21:E21").Copy() | Out-Null
$ws.Range("A22:E22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Step 2: re-style row 21 to match row 20's bordered look (style 10/11) ---
$ws.Range("A20:E20").Copy() | Out-Null
$ws.Range("A21:E21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

# --- Step 3: populate the new row 22 with the added entry ---
# (written in shared-string insertion order: English, filename, converted, translated)
$ws.Cells.Item(22, 3).Value2 = ' Mr. [CS:N]Drowzee[CR] told me to tell you,\n"Take care and thanks," [hero] and\n[partner].'
$ws.Cells.Item(22, 1).Value2 = 'SCRIPT/T01P02A/us2302.ssb'
$ws.Cells.Item(22, 5).Value2 = ' Íéòóåñ [CS:N]Äñïôèé[CR] ðïðñïòéì íåîÿ\nðåñåäàóû âàí åãï òìïâà: \"Áåñåãéóå òåáÿ é\nòðàòéáï âàí, [hero] é [partner]\".'
$ws.Cells.Item(22, 4).Value2 = ' Мистер [CS:N]Дроузи[CR] попросил меня\nпередать вам его слова: \"Берегите себя и\nспасибо вам, [hero] и [partner]\".'
$ws.Cells.Item(22, 2).Value2 = 19

# --- Step 4: match row height used throughout the sheet ---
$ws.Rows.Item(22).RowHeight = 43.2

# --- Step 5: scroll / selection bookkeeping, matching where Excel would land
#     after adding a row at the bottom of the table ---
$win = $excel.ActiveWindow
$win.ScrollRow = 21
$ws.Range("D25").Select() | Out-Null
